# Apply the reordering of species-record rows 47-50.
# The values in columns A, B, E, F, G, H, Q, R for these four rows
# have been cyclically/reverse-rotated: row47<->row50, row48<->row49.
# All other columns in these rows stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Target values after the edit (taken directly from the diff)
$data = @{
    47 = @{ A = 111974185; B = 90660; E = 4362; F = "Blå taggsvamp";      G = "Hydnellum caeruleum";  H = "(Hornem.) P.Karst.";                         Q = 439827.4842555065; R = 6952232.676732311 }
    48 = @{ A = 111974186; B = 90682; E = 2059; F = "Skrovlig taggsvamp"; G = "Hydnellum scabrosum";  H = "(Fr.) E.Larss., K.H.Larss. & Kõljalg";       Q = 439860.448822267;  R = 6952249.98427855 }
    49 = @{ A = 111974191; B = 90652; E = 3100; F = "Talltaggsvamp";      G = "Bankera fuligineoalba"; H = "(Schmidt : Fr.) Pouzar";                     Q = 439977.5118376439; R = 6952213.872195411 }
    50 = @{ A = 111974188; B = 90652; E = 3100; F = "Talltaggsvamp";      G = "Bankera fuligineoalba"; H = "(Schmidt : Fr.) Pouzar";                     Q = 439869.6589509377; R = 6952225.479112641 }
}

foreach ($row in 47..50) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}
